$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-09-08 Monday" "2025-09-09 Tuesday"

Replace-Text "61×48=" "96×42="
Replace-Text "89×30=" "87×56="
Replace-Text "49×48=" "61×22="
Replace-Text "82×98=" "28×90="
Replace-Text "31×60=" "35×11="
Replace-Text "59×19=" "59×43="
Replace-Text "65×57=" "95×53="
Replace-Text "22×19=" "28×49="
Replace-Text "97×41=" "17×92="
Replace-Text "21×86=" "74×61="
Replace-Text "87×68=" "84×32="
Replace-Text "51×85=" "56×17="
Replace-Text "76×77=" "73×73="
Replace-Text "73×40=" "48×21="
Replace-Text "70×48=" "31×31="
Replace-Text "39×65=" "75×30="
Replace-Text "89×55=" "35×70="
Replace-Text "56×33=" "27×50="
Replace-Text "81×87=" "65×78="
Replace-Text "41×81=" "46×62="
Replace-Text "21×52=" "80×68="
Replace-Text "73×44=" "52×52="
Replace-Text "66×18=" "68×25="
Replace-Text "65×16=" "13×76="
Replace-Text "48×95=" "65×49="
